$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column before column A to make room for "code"
$ws.Columns("A:A").Insert()

# Set header for new column A
$ws.Range("A1").Value = "code"
$ws.Range("A2").Value = 11
$ws.Range("A3").Value = 12

# Update selection to A4
$ws.Range("A4").Select()
